$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new column J: copy header style/format from I2, set header text "Increment"
$ws.Range("I2").Copy()
$ws.Range("J2").PasteSpecial(-4122)
$ws.Range("J2").Value2 = "Increment"

# Copy style/format from B3 (s=3) for J3 data cell, set value "1/1"
$ws.Range("B3").Copy()
$ws.Range("J3").PasteSpecial(-4122)
$ws.Range("J3").Value2 = "1/1"

# Update title text last, so new shared strings are appended in the right order
$ws.Range("A1").Value2 = "CLIENT Logo & Info"

# Extend the title band formatting and merged range to include column J
$ws.Range("I1").Copy()
$ws.Range("J1").PasteSpecial(-4122)
$ws.Range("A1:I1").UnMerge()
$ws.Range("A1:J1").Merge()

# Size the new column similarly to its neighbours
$ws.Columns("J").ColumnWidth = 11.95

# Update the active selection to match the template
$ws.Range("F4").Select() | Out-Null
